$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 420 entirely ("「我らも米を食います」" post), shifting all
# subsequent rows up by one.
$ws.Rows("420:420").Delete()
